$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AverageAccuracy")
$ws.Range("A1").Value = 0.76744186046511631
